$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "32+19="  # was 82-59=
$t.Cell(1,2).Range.Text = "8+91="  # was 8+73=
$t.Cell(1,3).Range.Text = "1+51="  # was 70-42=
$t.Cell(1,4).Range.Text = "93-9="  # was 28+63=
$t.Cell(1,5).Range.Text = "91+3="  # was 68-16=
$t.Cell(2,1).Range.Text = "86-44="  # was 11+6=
$t.Cell(2,2).Range.Text = "29+18="  # was 28-3=
$t.Cell(2,3).Range.Text = "43-0="  # was 9+15=
$t.Cell(2,4).Range.Text = "48-21="  # was 5+51=
$t.Cell(2,5).Range.Text = "76-19="  # was 20+1=
$t.Cell(3,1).Range.Text = "89-62="  # was 71+26=
$t.Cell(3,2).Range.Text = "39-26="  # was 29+38=
$t.Cell(3,3).Range.Text = "17+54="  # was 59-45=
$t.Cell(3,4).Range.Text = "79-67="  # was 38+47=
$t.Cell(3,5).Range.Text = "81-51="  # was 4+57=
$t.Cell(4,1).Range.Text = "9+39="  # was 12+38=
$t.Cell(4,2).Range.Text = "37+56="  # was 81-4=
$t.Cell(4,3).Range.Text = "0+56="  # was 68+9=
$t.Cell(4,4).Range.Text = "8+21="  # was 74-36=
$t.Cell(4,5).Range.Text = "80+5="  # was 54-25=
$t.Cell(5,1).Range.Text = "17+59="  # was 85-33=
$t.Cell(5,2).Range.Text = "83-66="  # was 66-56=
$t.Cell(5,3).Range.Text = "62+15="  # was 34+43=
$t.Cell(5,4).Range.Text = "1+74="  # was 73-63=
$t.Cell(5,5).Range.Text = "55+8="  # was 13+13=
$t.Cell(6,1).Range.Text = "26+20="  # was 22+29=
$t.Cell(6,2).Range.Text = "61-29="  # was 17+76=
$t.Cell(6,3).Range.Text = "86-85="  # was 3+2=
$t.Cell(6,4).Range.Text = "68-35="  # was 57-48=
$t.Cell(6,5).Range.Text = "80-46="  # was 67-6=
$t.Cell(7,1).Range.Text = "49+31="  # was 8+64=
$t.Cell(7,2).Range.Text = "51-40="  # was 0+11=
$t.Cell(7,3).Range.Text = "36-1="  # was 68+19=
$t.Cell(7,4).Range.Text = "8-6="  # was 21+26=
$t.Cell(7,5).Range.Text = "51-6="  # was 91-11=
$t.Cell(8,1).Range.Text = "20+76="  # was 78-45=
$t.Cell(8,2).Range.Text = "4+7="  # was 1+24=
$t.Cell(8,3).Range.Text = "67+7="  # was 37+37=
$t.Cell(8,4).Range.Text = "1+91="  # was 14+44=
$t.Cell(8,5).Range.Text = "57+29="  # was 3+74=
$t.Cell(9,1).Range.Text = "74-70="  # was 79-36=
$t.Cell(9,2).Range.Text = "15+60="  # was 36+37=
$t.Cell(9,3).Range.Text = "4-2="  # was 30+26=
$t.Cell(9,4).Range.Text = "99-6="  # was 35+61=
$t.Cell(9,5).Range.Text = "51+21="  # was 64-48=
$t.Cell(10,1).Range.Text = "47+17="  # was 37+7=
$t.Cell(10,2).Range.Text = "37+45="  # was 13+11=
$t.Cell(10,3).Range.Text = "6+68="  # was 82+7=
$t.Cell(10,4).Range.Text = "87-66="  # was 10+55=
$t.Cell(10,5).Range.Text = "27-24="  # was 2+18=
$t.Cell(11,1).Range.Text = "83-65="  # was 30-20=
$t.Cell(11,2).Range.Text = "45+53="  # was 0+6=
$t.Cell(11,3).Range.Text = "71+18="  # was 89-44=
$t.Cell(11,4).Range.Text = "63-12="  # was 60-56=
$t.Cell(11,5).Range.Text = "25+63="  # was 96-43=
$t.Cell(12,1).Range.Text = "47+21="  # was 35+3=
$t.Cell(12,2).Range.Text = "37+30="  # was 6+4=
$t.Cell(12,3).Range.Text = "15-2="  # was 86-57=
$t.Cell(12,4).Range.Text = "80-33="  # was 57-10=
$t.Cell(12,5).Range.Text = "92-35="  # was 13+2=
$t.Cell(13,1).Range.Text = "6+56="  # was 68-9=
$t.Cell(13,2).Range.Text = "20+36="  # was 85-77=
$t.Cell(13,3).Range.Text = "67+5="  # was 36+59=
$t.Cell(13,4).Range.Text = "32-1="  # was 2+39=
$t.Cell(13,5).Range.Text = "33+38="  # was 94-2=
$t.Cell(14,1).Range.Text = "13+47="  # was 61-11=
$t.Cell(14,2).Range.Text = "22+52="  # was 53+35=
$t.Cell(14,3).Range.Text = "4+81="  # was 89-37=
$t.Cell(14,4).Range.Text = "29-25="  # was 34+0=
$t.Cell(14,5).Range.Text = "46+21="  # was 30+12=
$t.Cell(15,1).Range.Text = "92+7="  # was 18+18=
$t.Cell(15,2).Range.Text = "57-28="  # was 89-79=
$t.Cell(15,3).Range.Text = "88-51="  # was 35+49=
$t.Cell(15,4).Range.Text = "54-41="  # was 82+13=
$t.Cell(15,5).Range.Text = "81-35="  # was 49+39=
$t.Cell(16,1).Range.Text = "20+52="  # was 27+4=
$t.Cell(16,2).Range.Text = "0+78="  # was 65-4=
$t.Cell(16,3).Range.Text = "88-26="  # was 71+27=
$t.Cell(16,4).Range.Text = "8-4="  # was 34+62=
$t.Cell(16,5).Range.Text = "48+2="  # was 97-27=
$t.Cell(17,1).Range.Text = "8+16="  # was 47-15=
$t.Cell(17,2).Range.Text = "49-7="  # was 34+53=
$t.Cell(17,3).Range.Text = "13+70="  # was 34-7=
$t.Cell(17,4).Range.Text = "99-91="  # was 61-5=
$t.Cell(17,5).Range.Text = "78-74="  # was 43-27=
$t.Cell(18,1).Range.Text = "54+42="  # was 75-11=
$t.Cell(18,2).Range.Text = "52+0="  # was 20+30=
$t.Cell(18,3).Range.Text = "5+40="  # was 43+16=
$t.Cell(18,4).Range.Text = "85-65="  # was 98-15=
$t.Cell(18,5).Range.Text = "37-36="  # was 97-85=
$t.Cell(19,1).Range.Text = "18-2="  # was 68-61=
$t.Cell(19,2).Range.Text = "27+34="  # was 74-5=
$t.Cell(19,3).Range.Text = "27-21="  # was 32+24=
$t.Cell(19,4).Range.Text = "65-33="  # was 89+0=
$t.Cell(19,5).Range.Text = "59+9="  # was 29+40=
$t.Cell(20,1).Range.Text = "23+23="  # was 24-22=
$t.Cell(20,2).Range.Text = "47-28="  # was 43+3=
$t.Cell(20,3).Range.Text = "86-36="  # was 95-12=
$t.Cell(20,4).Range.Text = "5+62="  # was 71-56=
$t.Cell(20,5).Range.Text = "26-12="  # was 59-6=
